$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two species records in rows 2 and 3 were swapped (their identity /
# taxon / Ost / Nord fields), while the shared site-level fields (C, D, I,
# K, P, S, T, U, V, W, Y, AD, AE, AG, AT, AW, AX, AY) stay where they are.
# Row 2 becomes what used to be row 3, and vice versa. The Ost/Nord
# coordinates are also rounded to whole numbers as part of this update, and
# the Starttid/Sluttid (Z/AB) values are removed entirely.

# --- Row 2 (now holds what used to be row 3's record) ---
$ws.Range("A2").Value = 111950173
$ws.Range("B2").Value = 90658
$ws.Range("E2").Value = 4361
$ws.Range("F2").Value = "Orange taggsvamp"
$ws.Range("G2").Value = "Hydnellum aurantiacum"
$ws.Range("H2").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q2").Value = 465440
$ws.Range("R2").Value = 6875680

# --- Row 3 (now holds what used to be row 2's record) ---
$ws.Range("A3").Value = 111950243
$ws.Range("B3").Value = 90689
$ws.Range("E3").Value = 5966
$ws.Range("F3").Value = "Motaggsvamp"
$ws.Range("G3").Value = "Sarcodon squamosus"
$ws.Range("H3").Value = "(Schaeff.) Quél."
$ws.Range("Q3").Value = 465473
$ws.Range("R3").Value = 6875785

# --- Remove the Starttid (Z) and Sluttid (AB) values on both rows ---
$ws.Range("Z2:Z3").ClearContents()
$ws.Range("AB2:AB3").ClearContents()
